$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" data column (Q) to the table, mirroring the formatting
# already used by the existing year columns. Missing/unavailable data points
# are written as "-" (matching the convention already used elsewhere in the
# sheet for absent values).

# Q4: year header, same look as the other year headers (P4 = 2019)
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Q5: first data row (country total), same look as P5
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 0.1

# Q6: Batken oblast
$ws.Range("D6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 0.2

# Q7: Djalal-Abad oblast - value not available
$ws.Range("D6").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = "-"

# Q8: Issyk-Kul oblast
$ws.Range("D6").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 0.2

# Q9: Naryn oblast - value not available
$ws.Range("D6").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = "-"

# Q10: Osh oblast
$ws.Range("D6").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 0.1

# Q11: Talas oblast - value not available
$ws.Range("D6").Copy()
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q11").Value = "-"

# Q12: Chui oblast
$ws.Range("D6").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = 0.3

# Q13: Bishkek city - value not available, same look as P13
$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q13").Value = "-"

# Q14: Osh city - value not available, same look as P14
$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = "-"

$excel.CutCopyMode = $false

$ws.Range("O17").Select()
